$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are assisting an international trade attorney in settling a dispute regarding importing poultry meat from Brazil to the European Union. The Brazilian government regarded the EU's import criteria as inconsistent with applicable international trade laws, and the dispute is in its initial stages. One of your clients urgently needs to know what consequences the dispute will have in the short term because they produce poultry meat.What do you include in your answer to the client?",
        "ques_type": 2,
        "options": [
            "The dispute will not affect trade for the time being because the first step in dispute settlement is the initiation of consultations at the World Trade Organization (WTO).",
            "The WTO will decide in favor of Brazil because the organization is enforcing free trade.",
            "The WTO will decide in favor of the EU because EU trade law is not bound to WTO free trade law. ",
            "Brazilian producers will move production to the EU to bypass restrictions."
        ],
        "score": "The dispute will not affect trade for the time being because the first step in dispute settlement is the initiation of consultations at the World Trade Organization (WTO)."
    },
    {
        "title": "You are a legal consultant at a German car manufacturer in the US, and the US government is planning to introduce new legal requirements for diesel cars. Your manager asks you to gather information and initiate lobbying activities to minimize the risk of your products being negatively affected by the planned legal reform. Which of the following organizations will you approach to help you achieve your goal?",
        "ques_type": 2,
        "options": [
            "The US embassy in Germany.",
            "The German embassy in the US.",
            "Delegation of German Industry and Commerce in Washington, DC.",
            "The International Trade Administration of the US government in Germany."
        ],
        "score": "Delegation of German Industry and Commerce in Washington, DC."
    },
    {
        "title": "You are an international trade policy analyst for a US brewer. Your brand has successfully exported beer to the European Union market and is gaining popularity in the entire EU. The German Brewers Association has initiated legal action against your product, arguing it does not meet the legal criteria to be categorized as \u201cbeer\u201d in Germany. Which line of argument will you apply to avoid negative legal consequences for your product?",
        "ques_type": 2,
        "options": [
            "The German Brewers Association has the right to restrict the sale of certain products.",
            "German law is exclusively bound to EU law, and the association's position is irrelevant. ",
            "Germany and the US are members of the World Trade Organization (WTO), and thus the WTO is responsible for settling the dispute. ",
            "The case falls under the National Treatment principle, making the association's position irrelevant. "
        ],
        "score": "The case falls under the National Treatment principle, making the association's position irrelevant."
    },
    {
        "title": "You are an international trade law associate for a corn producer in Romania. Your company intends to export its corn to Italy. The Italian national bodies notify your company that it deems your corn type to be incompatible with import regulations even though the US has been exporting the identical strain to Italy for years. Which line of argument do you apply to avoid exclusion from the Italian market?",
        "ques_type": 2,
        "options": [
            "Italy has the right to limit agricultural imports regardless of other import arrangements to protect its farmers and agricultural market. ",
            "The Romanian producer should consider growing other corn types to export to Italy. ",
            "Your company should look into exporting to a third-party country and re-export its product from there. ",
            "The Italian claim is inconsistent with the Most Favored Nation principle, and your company has the right to export to the Italian market. "
        ],
        "score": "The Italian claim is inconsistent with the Most Favored Nation principle, and your company has the right to export to the Italian market."
    }
]
'@

# Remove the old row 1 (A1 = 0 with bold/border style) entirely, shifting A2 up to A1
$ws.Rows.Item(1).Delete() | Out-Null

# Set A1 to the updated shared string text (row 2's former content, now in row 1 with default style)
$ws.Range("A1").Value = $text

# The long multi-line text auto-expands the row height; restore it to the standard height
$ws.Rows.Item(1).AutoFit() | Out-Null
